$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the two section-header rows ("grandes regiões e unidades da federação" at
# row 8, "situação do domicílio" at row 5). Delete the higher-numbered row first so
# the lower row's index doesn't shift before it is removed.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
